$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text / date-as-text fields
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("M2").Value = "2020-12-15 00:00:00"
$ws.Range("N2").Value = "2019-12-31 00:00:00"

# Numeric fields
$ws.Range("O2").Value = 1390618711.43
$ws.Range("P2").Value = 479097142.71
$ws.Range("Q2").Value = 287832012.28
$ws.Range("R2").Value = -6.5958583509
$ws.Range("S2").Value = 37302689.14
$ws.Range("T2").Value = -8.317492528700001
$ws.Range("U2").Value = 59601198.76
$ws.Range("V2").Value = 22.7981583085
$ws.Range("W2").Value = 729391693.34
$ws.Range("X2").Value = 176697912.98
$ws.Range("Y2").Value = 17.9426022997
$ws.Range("Z2").Value = 220864640.47
$ws.Range("AA2").Value = 12.6457557269
$ws.Range("AB2").Value = 661227018.09
$ws.Range("AC2").Value = 35.6392482557
$ws.Range("AD2").Value = 36.8922604734
$ws.Range("AE2").Value = 38.0483485573
$ws.Range("AF2").Value = 79.4645872948
$ws.Range("AG2").Value = 52.4508758112
